# Rename the sheets: "Men's" -> "MALE", "Women's" -> "FEMALE"
$wb = $excel.ActiveWorkbook
$wsMale = $wb.Worksheets.Item(1)
$wsMale.Name = "MALE"
$wsFemale = $wb.Worksheets.Item(2)
$wsFemale.Name = "FEMALE"

# Update the selection on the FEMALE sheet (was B4, now C33) before switching
# the active tab away from it, so the new selection is persisted.
$wsFemale.Range("C33").Select()

# Make the MALE sheet the active / selected tab (previously FEMALE was active).
$wsMale.Activate()
